$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 88.5
$ws.Range("I11").Value = 88.5
$ws.Range("K11").Value = 88.5
$ws.Range("M11").Value = 51.5
$ws.Range("H96").Value = 1570.3684
$ws.Range("I96").Value = 883.8
$ws.Range("K96").Value = 2651.4
$ws.Range("M96").Value = -1278.4
$ws.Range("H99").Value = 824.1429000000001
$ws.Range("I99").Value = 675.8
$ws.Range("J99").Value = 1195
$ws.Range("K99").Value = 2027.4
$ws.Range("L99").Value = 3585
$ws.Range("M99").Value = -529.3999999999999
$ws.Range("N99").Value = -6581
$ws.Range("H107").Value = 5632.4614
$ws.Range("I107").Value = 5458.609
$ws.Range("K107").Value = 5458.609
$ws.Range("M107").Value = -3538.609
$ws.Range("H111").Value = 1274.2222
$ws.Range("I111").Value = 1343.2
$ws.Range("K111").Value = 4029.6
$ws.Range("M111").Value = -962.6000000000004
$ws.Range("H118").Value = 1156.25
$ws.Range("I118").Value = 575
$ws.Range("J118").Value = 2900
$ws.Range("K118").Value = 1725
$ws.Range("L118").Value = 8700
$ws.Range("M118").Value = -68
$ws.Range("N118").Value = -12014
$ws.Range("H125").Value = 1043.4
$ws.Range("I125").Value = 1102.6471
$ws.Range("K125").Value = 9923.823899999999
$ws.Range("M125").Value = -7463.823899999999
$ws.Range("H131").Value = 8105.467
$ws.Range("I131").Value = 3038.476
$ws.Range("J131").Value = 19928.445
$ws.Range("K131").Value = 9115.428
$ws.Range("L131").Value = 59785.335
$ws.Range("M131").Value = -4075.428
$ws.Range("N131").Value = -69865.33499999999
$ws.Range("H132").Value = 1829.6818
$ws.Range("I132").Value = 1842.1428
$ws.Range("J132").Value = 1759.9
$ws.Range("K132").Value = 5526.428400000001
$ws.Range("L132").Value = 5279.700000000001
$ws.Range("M132").Value = -2996.428400000001
$ws.Range("N132").Value = -10339.7
$ws.Range("H138").Value = 2859.5088
$ws.Range("J138").Value = 4180.0967
$ws.Range("L138").Value = 12540.2901
$ws.Range("N138").Value = -22820.2901

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 102.23529
$ws.Range("I5").Value = 93.166664
$ws.Range("K5").Value = 93.166664
$ws.Range("M5").Value = 18.833336
$ws.Range("H14").Value = 10814.667
$ws.Range("I14").Value = 11977.6
$ws.Range("K14").Value = 11977.6
$ws.Range("M14").Value = -11802.6
$ws.Range("H110").Value = 3425.7932
$ws.Range("I110").Value = 3190.3845
$ws.Range("J110").Value = 5466
$ws.Range("K110").Value = 3190.3845
$ws.Range("L110").Value = 5466
$ws.Range("M110").Value = -1145.3845
$ws.Range("N110").Value = -9556

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 102.23529
$ws.Range("I4").Value = 93.166664
$ws.Range("K4").Value = 93.166664
$ws.Range("M4").Value = 21.833336

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10580
$ws.Range("I4").Value = 11000
$ws.Range("J4").Value = 10475
$ws.Range("K4").Value = 11000
$ws.Range("L4").Value = 10475
$ws.Range("M4").Value = -10888
$ws.Range("N4").Value = -10699
$ws.Range("H21").Value = 9974
$ws.Range("J21").Value = 9974
$ws.Range("L21").Value = 9974
$ws.Range("N21").Value = -10444
$ws.Range("H25").Value = 3000.7058
$ws.Range("I25").Value = 1458.7858
$ws.Range("J25").Value = 10196.333
$ws.Range("K25").Value = 1458.7858
$ws.Range("L25").Value = 10196.333
$ws.Range("M25").Value = -1284.7858
$ws.Range("N25").Value = -10544.333
$ws.Range("H58").Value = 1840.4
$ws.Range("I58").Value = 1050.75
$ws.Range("K58").Value = 1050.75
$ws.Range("M58").Value = -847.75
$ws.Range("H105").Value = 2098.6
$ws.Range("I105").Value = 2098.6
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2098.6
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -351.5999999999999
$ws.Range("N105").ClearContents() | Out-Null
$ws.Range("H136").Value = 1840.4
$ws.Range("I136").Value = 1050.75
$ws.Range("K136").Value = 3152.25
$ws.Range("M136").Value = -602.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 17363354
$ws.Range("I4").Value = 134570.38
$ws.Range("K4").Value = 403711.14
$ws.Range("M4").Value = -403599.14
$ws.Range("H12").Value = 706.5161000000001
$ws.Range("J12").Value = 790.1111
$ws.Range("L12").Value = 2370.3333
$ws.Range("N12").Value = -2716.3333
$ws.Range("H26").Value = 2520.5
$ws.Range("I26").Value = 2520.5
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 7561.5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -7273.5
$ws.Range("N26").ClearContents() | Out-Null
$ws.Range("H36").Value = 4589.3
$ws.Range("I36").Value = 895.4286
$ws.Range("J36").Value = 6578.3076
$ws.Range("K36").Value = 2686.2858
$ws.Range("L36").Value = 19734.9228
$ws.Range("M36").Value = -2517.2858
$ws.Range("N36").Value = -20072.9228
$ws.Range("H80").Value = 11001.35
$ws.Range("I80").Value = 16456.25
$ws.Range("J80").Value = 7364.75
$ws.Range("K80").Value = 49368.75
$ws.Range("L80").Value = 22094.25
$ws.Range("M80").Value = -48432.75
$ws.Range("N80").Value = -23966.25
$ws.Range("H83").Value = 11001.35
$ws.Range("I83").Value = 16456.25
$ws.Range("J83").Value = 7364.75
$ws.Range("K83").Value = 148106.25
$ws.Range("L83").Value = 66282.75
$ws.Range("M83").Value = -143426.25
$ws.Range("N83").Value = -75642.75
$ws.Range("H131").Value = 1835.0448
$ws.Range("J131").Value = 2028.4728
$ws.Range("L131").Value = 6085.4184
$ws.Range("N131").Value = -16165.4184

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 14510.8
$ws.Range("I18").Value = 14510.8
$ws.Range("K18").Value = 14510.8
$ws.Range("M18").Value = -14217.8
$ws.Range("H102").Value = 4008.796
$ws.Range("I102").Value = 3710.0227
$ws.Range("K102").Value = 3710.0227
$ws.Range("M102").Value = -2088.0227
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents() | Out-Null
$ws.Range("H126").Value = 5004.8
$ws.Range("I126").Value = 4506
$ws.Range("J126").Value = 7000
$ws.Range("K126").Value = 13518
$ws.Range("L126").Value = 21000
$ws.Range("M126").Value = -11048
$ws.Range("N126").Value = -25940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 28135898
$ws.Range("I2").Value = 150002460
$ws.Range("K2").Value = 150002460
$ws.Range("M2").Value = -150002348
$ws.Range("H7").Value = 4311.4
$ws.Range("I7").Value = 4311.4
$ws.Range("K7").Value = 4311.4
$ws.Range("M7").Value = -4199.4
$ws.Range("H23").Value = 20000
$ws.Range("I23").Value = 20000
$ws.Range("K23").Value = 20000
$ws.Range("M23").Value = -19770
$ws.Range("H61").Value = 13772.308
$ws.Range("I61").Value = 11512.909
$ws.Range("J61").Value = 26199
$ws.Range("K61").Value = 11512.909
$ws.Range("L61").Value = 26199
$ws.Range("M61").Value = -11310.909
$ws.Range("N61").Value = -26603
$ws.Range("H113").Value = 13772.308
$ws.Range("I113").Value = 11512.909
$ws.Range("J113").Value = 26199
$ws.Range("K113").Value = 11512.909
$ws.Range("L113").Value = 26199
$ws.Range("M113").Value = -9342.909
$ws.Range("N113").Value = -30539
$ws.Range("H126").Value = 4311.4
$ws.Range("I126").Value = 4311.4
$ws.Range("K126").Value = 12934.2
$ws.Range("M126").Value = -10464.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 4446.533
$ws.Range("I14").Value = 3764.1785
$ws.Range("J14").Value = 13999.5
$ws.Range("K14").Value = 3764.1785
$ws.Range("L14").Value = 13999.5
$ws.Range("M14").Value = -3596.1785
$ws.Range("N14").Value = -14335.5
$ws.Range("H15").Value = 18000
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents() | Out-Null
$ws.Range("H100").Value = 1004.94446
$ws.Range("I100").Value = 819.3333
$ws.Range("K100").Value = 1638.6666
$ws.Range("M100").Value = -1097.6666
$ws.Range("H126").Value = 1835529.8
$ws.Range("I126").Value = 2384678.2
$ws.Range("J126").Value = 5034.6665
$ws.Range("K126").Value = 7154034.600000001
$ws.Range("L126").Value = 15103.9995
$ws.Range("M126").Value = -7151564.600000001
$ws.Range("N126").Value = -20043.9995
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents() | Out-Null
$ws.Range("N132").ClearContents() | Out-Null
